# Update leve-profit calculation columns (H:N) across several sheets
# to reflect refreshed market-board pricing data, per scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1455.32
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1580.1428
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4740.428400000001
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6956.428400000001

# Row 137
$ws.Range("H137").Value = 2854.077
$ws.Range("I137").Value = 3410.9285
$ws.Range("J137").Value = 2204.4167
$ws.Range("K137").Value = 10232.7855
$ws.Range("L137").Value = 6613.250100000001
$ws.Range("M137").Value = -7682.7855
$ws.Range("N137").Value = -11713.2501

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4188.409
$ws.Range("I61").Value = 3089.611
$ws.Range("J61").Value = 9133
$ws.Range("K61").Value = 3089.611
$ws.Range("L61").Value = 9133
$ws.Range("M61").Value = -2877.611
$ws.Range("N61").Value = -9557

# Row 74
$ws.Range("H74").Value = 1425.24
$ws.Range("I74").Value = 1205.973
$ws.Range("J74").Value = 2049.3076
$ws.Range("K74").Value = 1205.973
$ws.Range("L74").Value = 2049.3076
$ws.Range("M74").Value = -331.973
$ws.Range("N74").Value = -3797.3076

# Row 77
$ws.Range("H77").Value = 1425.24
$ws.Range("I77").Value = 1205.973
$ws.Range("J77").Value = 2049.3076
$ws.Range("K77").Value = 6029.865
$ws.Range("L77").Value = 10246.538
$ws.Range("M77").Value = -1661.865
$ws.Range("N77").Value = -18982.538

# Row 132
$ws.Range("H132").Value = 2645.4375
$ws.Range("I132").Value = 1154.1052
$ws.Range("J132").Value = 4825.077
$ws.Range("K132").Value = 3462.3156
$ws.Range("L132").Value = 14475.231
$ws.Range("M132").Value = -932.3155999999999
$ws.Range("N132").Value = -19535.231

# Row 136
$ws.Range("H136").Value = 4188.409
$ws.Range("I136").Value = 3089.611
$ws.Range("J136").Value = 9133
$ws.Range("K136").Value = 9268.832999999999
$ws.Range("L136").Value = 27399
$ws.Range("M136").Value = -6718.832999999999
$ws.Range("N136").Value = -32499

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 806.9091
$ws.Range("I107").Value = 764
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 764
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1156
$ws.Range("N107").Value = -4840

# Row 134
$ws.Range("H134").Value = 2407.2974
$ws.Range("I134").Value = 1608.6428
$ws.Range("J134").Value = 4892
$ws.Range("K134").Value = 4825.928400000001
$ws.Range("L134").Value = 14676
$ws.Range("M134").Value = -2290.928400000001
$ws.Range("N134").Value = -19746

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11254.103
$ws.Range("I31").Value = 1192.4348
$ws.Range("J31").Value = 25717.75
$ws.Range("K31").Value = 1192.4348
$ws.Range("L31").Value = 25717.75
$ws.Range("M31").Value = -897.4348
$ws.Range("N31").Value = -26307.75

# Row 34
$ws.Range("H34").Value = 11254.103
$ws.Range("I34").Value = 1192.4348
$ws.Range("J34").Value = 25717.75
$ws.Range("K34").Value = 1192.4348
$ws.Range("L34").Value = 25717.75
$ws.Range("M34").Value = -990.4348
$ws.Range("N34").Value = -26121.75

# Row 106
$ws.Range("H106").Value = 35888.75
$ws.Range("I106").Value = 20000
$ws.Range("J106").Value = 41185
$ws.Range("K106").Value = 20000
$ws.Range("L106").Value = 41185
$ws.Range("M106").Value = -18738
$ws.Range("N106").Value = -43709

# Row 134
$ws.Range("H134").Value = 7861.8
$ws.Range("I134").Value = 8800.666999999999
$ws.Range("K134").Value = 26402.001
$ws.Range("M134").Value = -23867.001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2054.1875
$ws.Range("I5").Value = 1002
$ws.Range("J5").Value = 2404.9167
$ws.Range("K5").Value = 3006
$ws.Range("L5").Value = 7214.750100000001
$ws.Range("M5").Value = -2894
$ws.Range("N5").Value = -7438.750100000001

# Row 122
$ws.Range("H122").Value = 859.61536
$ws.Range("J122").Value = 1566.3334
$ws.Range("L122").Value = 14097.0006
$ws.Range("N122").Value = -18997.0006

# Row 125
$ws.Range("H125").Value = 1447.1428

# Row 131
$ws.Range("H131").Value = 1319.8628
$ws.Range("J131").Value = 1424.3864
$ws.Range("L131").Value = 4273.1592
$ws.Range("N131").Value = -14353.1592

# Row 135
$ws.Range("H135").Value = 2054.1875
$ws.Range("I135").Value = 1002
$ws.Range("J135").Value = 2404.9167
$ws.Range("K135").Value = 9018
$ws.Range("L135").Value = 21644.2503
$ws.Range("M135").Value = -6483
$ws.Range("N135").Value = -26714.2503

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3329.125
$ws.Range("I132").Value = 2493.5
$ws.Range("J132").Value = 3746.9375
$ws.Range("K132").Value = 7480.5
$ws.Range("L132").Value = 11240.8125
$ws.Range("M132").Value = -4950.5
$ws.Range("N132").Value = -16300.8125

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 17893866
$ws.Range("I46").Value = 55667376
$ws.Range("J46").Value = 1150.8948
$ws.Range("K46").Value = 55667376
$ws.Range("L46").Value = 1150.8948
$ws.Range("M46").Value = -55667188
$ws.Range("N46").Value = -1526.8948

# Row 61
$ws.Range("H61").Value = 1882.5
$ws.Range("I61").Value = 1488.5714
$ws.Range("J61").Value = 2801.6667
$ws.Range("K61").Value = 1488.5714
$ws.Range("L61").Value = 2801.6667
$ws.Range("M61").Value = -1286.5714
$ws.Range("N61").Value = -3205.6667

# Row 113
$ws.Range("H113").Value = 1882.5
$ws.Range("I113").Value = 1488.5714
$ws.Range("J113").Value = 2801.6667
$ws.Range("K113").Value = 1488.5714
$ws.Range("L113").Value = 2801.6667
$ws.Range("M113").Value = 681.4286
$ws.Range("N113").Value = -7141.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 555.13635
$ws.Range("I113").Value = 572.3570999999999
$ws.Range("J113").Value = 525
$ws.Range("K113").Value = 1717.0713
$ws.Range("L113").Value = 1575
$ws.Range("M113").Value = 452.9287000000002
$ws.Range("N113").Value = -5915

# Row 132
$ws.Range("H132").Value = 2058.8667
$ws.Range("I132").Value = 1450.8077
$ws.Range("J132").Value = 2890.9473
$ws.Range("K132").Value = 4352.4231
$ws.Range("L132").Value = 8672.841899999999
$ws.Range("M132").Value = -1822.4231
$ws.Range("N132").Value = -13732.8419

# Row 136
$ws.Range("H136").Value = 4755.378
$ws.Range("I136").Value = 982.4286
$ws.Range("J136").Value = 5450.3945
$ws.Range("K136").Value = 2947.2858
$ws.Range("L136").Value = 16351.1835
$ws.Range("M136").Value = -397.2857999999997
$ws.Range("N136").Value = -21451.1835
